# LoopingMacro.xlsx edit: add a new "localdb" command-type column/category
# to the hidden '#system' sheet, with its 6 functions, insert it in the
# alphabetically-correct spot (between "json" and "macro"), and register
# the corresponding entries everywhere the command-type catalogue is kept
# (the "target" list in column A, and the named ranges in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Make room for the new "localdb" category: insert a new column at N,
#    pushing the existing N..AC ("macro".."xml") one column to the right
#    (O..AD). This also lets defined names keep matching single columns.
# ---------------------------------------------------------------------
$ws.Range("N1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2) Populate the new column N with the "localdb" header and its six
#    functions (content taken from the authoring commit).
# ---------------------------------------------------------------------
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 3) The "target" list in column A enumerates every command-type name in
#    alphabetical order. Shift rows 14..29 ("macro".."xml") down by one
#    (to 15..30) to make room, then insert "localdb" at row 14 (it sorts
#    right after "json" and before "macro").
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $moved = $ws.Range("A$r").Value()
    $ws.Range("A$($r + 1)").Value = $moved
}
$ws.Range("A14").Value = "localdb"

# ---------------------------------------------------------------------
# 4) Update the workbook-level named ranges: everything from "macro"
#    onward (alphabetically) shifts one column to the right, and a new
#    "localdb" name is registered for the freshly inserted column N.
# ---------------------------------------------------------------------
$wb.Names.Item("macro").RefersTo      = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Item("mail").RefersTo       = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo     = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo        = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo      = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo      = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo        = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo      = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo        = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo       = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo     = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo        = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AD`$2:`$AD`$21"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
